# Apply the "Add files via upload" edit to the AGR_DEMAND sheet:
#  - insert two new rows above the existing "AGR_DEM_LIV_CAT_DAIRY" block
#    to add a new AGR_BIOG demand entry (with its unit row underneath)
#  - change the unit shown in D5 from "kanimals" to "PJ"
#  - leave the active selection on D7 (as last edited by the author)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AGR_DEMAND")

# --- 1. Insert two blank rows at row 6 (pushes everything below down by 2) ---
$ws.Rows.Item(6).Resize(2).Insert()

# --- 2. Fill in the new "AGR_BIOG" row (row 6) ---
$ws.Range("B6").Value = "AGR_BIOG"
$ws.Range("C6").Value = "PL"
$ws.Range("D6").Value = 0.5

# --- 3. Match formatting of the two new rows to the rest of the table ---
$newRows = $ws.Range("B6:K7")
$newRows.Style = "Normal 10 15 2"
$newRows.Interior.Pattern = -4142      # xlPatternNone - no fill
$newRows.Borders.LineStyle = 0         # xlLineStyleNone - no border
$newRows.HorizontalAlignment = -4108   # xlCenter
$newRows.VerticalAlignment = -4108     # xlCenter
$newRows.WrapText = $true
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75

# Row 7 gets the thick bottom border that marks the end of the sub-block
$ws.Range("B7:K7").Borders.Item(9).LineStyle = 1
$ws.Range("B7:K7").Borders.Item(9).Weight = 3

# --- 4. Update the unit label in D5 from "kanimals" to "PJ" ---
$ws.Range("D5").Value = "PJ"

# --- 5. Restore the author's last selection ---
$ws.Range("D7").Select()
